$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 25 de Marzo de 2020 a las 06:16"

# --- Estados Unidos (row 6): updated case counts ---
$ws.Range("B6").Value = 54881
$ws.Range("C6").Value = 58
$ws.Range("E6").Value = 53721

# --- Tailandia overtakes Chile in total cases and moves up the (sorted) table ---
# Rows 33-35 were Chile, Polonia, Tailandia; Tailandia's new total (934) now
# exceeds Chile's (922), so it is re-sorted above Chile/Polonia.
$ws.Range("A33").Value = "Tailandia"
$ws.Range("B33").Value = 934
$ws.Range("C33").Value = 107
$ws.Range("D33").Value = 52
$ws.Range("E33").Value = 878
$ws.Range("F33").Value = 7
$ws.Range("G33").Value = 0
$ws.Range("H33").Value = 4

$ws.Range("A34").Value = "Chile"
$ws.Range("B34").Value = 922
$ws.Range("C34").Value = 0
$ws.Range("D34").Value = 17
$ws.Range("E34").Value = 903
$ws.Range("F34").Value = 7
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = 2

$ws.Range("A35").Value = "Polonia"
$ws.Range("B35").Value = 901
$ws.Range("C35").Value = 0
$ws.Range("D35").Value = 1
$ws.Range("E35").Value = 890
$ws.Range("F35").Value = 3
$ws.Range("G35").Value = 0
$ws.Range("H35").Value = 10

# --- Argelia (row 64): updated active/recovered counts ---
$ws.Range("D64").Value = 65
$ws.Range("E64").Value = 180

# --- Camboya (row 95): updated case counts ---
$ws.Range("B95").Value = 93
$ws.Range("C95").Value = 2
$ws.Range("E95").Value = 89

# --- Nueva Caledonia overtakes Etiopia/Tanzania in total cases and moves up ---
# Rows 142-145 were Etiopia, Tanzania, Mongolia, Nueva Caledonia; Nueva
# Caledonia's new total (12) now ties/overtakes Etiopia's, so it is
# re-sorted above Etiopia/Tanzania/Mongolia.
$ws.Range("A142").Value = "Nueva Caledonia"
$ws.Range("B142").Value = 12
$ws.Range("C142").Value = 2
$ws.Range("D142").Value = 0
$ws.Range("E142").Value = 12
$ws.Range("F142").Value = 0
$ws.Range("G142").Value = 0
$ws.Range("H142").Value = 0

$ws.Range("A143").Value = "Etiopia"
$ws.Range("B143").Value = 12
$ws.Range("C143").Value = 0
$ws.Range("D143").Value = 0
$ws.Range("E143").Value = 12
$ws.Range("F143").Value = 0
$ws.Range("G143").Value = 0
$ws.Range("H143").Value = 0

$ws.Range("A144").Value = "Tanzania"
$ws.Range("B144").Value = 12
$ws.Range("C144").Value = 0
$ws.Range("D144").Value = 0
$ws.Range("E144").Value = 12
$ws.Range("F144").Value = 0
$ws.Range("G144").Value = 0
$ws.Range("H144").Value = 0

$ws.Range("A145").Value = "Mongolia"
$ws.Range("B145").Value = 10
$ws.Range("C145").Value = 0
$ws.Range("D145").Value = 0
$ws.Range("E145").Value = 10
$ws.Range("F145").Value = 0
$ws.Range("G145").Value = 0
$ws.Range("H145").Value = 0
